$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new progress entry ---
$ws.Range("A2:B2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 41118
$ws.Range("B3").Value = "FIxed SPI Slave Select issue. Need to insert a small delay between end of SPI transfer and rising edge of SS signal. Without the delay, SS rises 4.5us after end of transfer. This would seem to be in spec according to datasheet, but may be affected due to not having set up UC3C oscillator and clock tree correctly."

$ws.Range("C3").Value = "SPI-transfer-1.png"
$ws.Range("D3").Value = "SPI-transfer-2.png"
$ws.Range("E3").Value = "SPI-transfer-3-SS-Rises-Too-Early.png"
$ws.Range("F3").Value = "SPI-transfer-4-SS-Rises-Too-Early.png"

$ws.Hyperlinks.Add($ws.Range("C3"), "SPI-transfer-1.png")
$ws.Hyperlinks.Add($ws.Range("D3"), "SPI-transfer-2.png")
$ws.Hyperlinks.Add($ws.Range("E3"), "SPI-transfer-3-SS-Rises-Too-Early.png")
$ws.Hyperlinks.Add($ws.Range("F3"), "SPI-transfer-4-SS-Rises-Too-Early.png")
